$d = $word.ActiveDocument

$d.Content.Find.Execute("For every d_bits add a parity bit", $true, $false, $false, $false, $false, $true, 1, $false, "TESTREPL", 2)
